# Generate Report for Archive
# The localization status report was regenerated: the "Status" column for
# the rows that are currently mid-pipeline flips from "Ready for handoff"
# to "In Translation" on every sheet that surfaces it (the Overview roll-up
# columns for each target language, plus each language sheet's own Status
# column). Excel auto-fits those Status columns to the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns (E = zh-cn, F = de-de)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-language detail sheets: Status column (C)
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter than
# "Ready for handoff".
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
